$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1719896.1
$ws.Range("J17").Value = 1719896.1
$ws.Range("L17").Value = 5159688.300000001
$ws.Range("N17").Value = -5160024.300000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1748730.4
$ws.Range("J80").Value = 2958917
$ws.Range("L80").Value = 8876751
$ws.Range("N80").Value = -8878747

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1748730.4
$ws.Range("J83").Value = 2958917
$ws.Range("L83").Value = 26630253
$ws.Range("N83").Value = -26640237

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1974.5
$ws.Range("J86").Value = 1974.5
$ws.Range("L86").Value = 1974.5
$ws.Range("N86").Value = -4220.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1139.2
$ws.Range("J88").Value = 1188
$ws.Range("L88").Value = 1188
$ws.Range("N88").Value = -2000

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1974.5
$ws.Range("J89").Value = 1974.5
$ws.Range("L89").Value = 9872.5
$ws.Range("N89").Value = -21104.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1139.2
$ws.Range("J91").Value = 1188
$ws.Range("L91").Value = 1188
$ws.Range("N91").Value = -3996

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 10057
$ws.Range("I106").Value = 12879.8
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 12879.8
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -12248.8
$ws.Range("N106").Value = -4262

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2980.5593
$ws.Range("I116").Value = 2654.3455
$ws.Range("K116").Value = 2654.3455
$ws.Range("M116").Value = 787.6545000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 5595
$ws.Range("I135").Value = 4432.467
$ws.Range("K135").Value = 39892.20299999999
$ws.Range("M135").Value = -37357.20299999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4553.864
$ws.Range("I141").Value = 3959.25
$ws.Range("K141").Value = 11877.75
$ws.Range("M141").Value = -6697.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10433.19
$ws.Range("I32").Value = 4160.08
$ws.Range("K32").Value = 4160.08
$ws.Range("M32").Value = -3873.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5266.2856
$ws.Range("I45").Value = 4976.1177
$ws.Range("K45").Value = 4976.1177
$ws.Range("M45").Value = -4599.1177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14702.827
$ws.Range("I61").Value = 7991.5293
$ws.Range("J61").Value = 24210.5
$ws.Range("K61").Value = 7991.5293
$ws.Range("L61").Value = 24210.5
$ws.Range("M61").Value = -7779.5293
$ws.Range("N61").Value = -24634.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H108").Value = 50396.727
$ws.Range("J108").Value = 50396.727
$ws.Range("L108").Value = 50396.727
$ws.Range("N108").Value = -58076.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 6404.8335
$ws.Range("I110").Value = 6685.8
$ws.Range("K110").Value = 6685.8
$ws.Range("M110").Value = -4640.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 14702.827
$ws.Range("I136").Value = 7991.5293
$ws.Range("J136").Value = 24210.5
$ws.Range("K136").Value = 23974.5879
$ws.Range("L136").Value = 72631.5
$ws.Range("M136").Value = -21424.5879
$ws.Range("N136").Value = -77731.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 26793.172
$ws.Range("I20").Value = 7799.4375
$ws.Range("J20").Value = 42787.895
$ws.Range("K20").Value = 7799.4375
$ws.Range("L20").Value = 42787.895
$ws.Range("M20").Value = -7552.4375
$ws.Range("N20").Value = -43281.895

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3149.76
$ws.Range("I94").Value = 1247.1765
$ws.Range("K94").Value = 1247.1765
$ws.Range("M94").Value = -796.1765

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3578.3
$ws.Range("I105").Value = 3578.3
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3578.3
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1831.3
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3769.6538
$ws.Range("I86").Value = 3192.1538
$ws.Range("J86").Value = 4347.154
$ws.Range("K86").Value = 3192.1538
$ws.Range("L86").Value = 4347.154
$ws.Range("M86").Value = -2069.1538
$ws.Range("N86").Value = -6593.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3769.6538
$ws.Range("I89").Value = 3192.1538
$ws.Range("J89").Value = 4347.154
$ws.Range("K89").Value = 15960.769
$ws.Range("L89").Value = 21735.77
$ws.Range("M89").Value = -10344.769
$ws.Range("N89").Value = -32967.77

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8623
$ws.Range("I99").Value = 2019.8
$ws.Range("J99").Value = 12750
$ws.Range("K99").Value = 2019.8
$ws.Range("L99").Value = 12750
$ws.Range("M99").Value = -521.8
$ws.Range("N99").Value = -15746

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 38873.184
$ws.Range("J108").Value = 38873.184
$ws.Range("L108").Value = 38873.184
$ws.Range("N108").Value = -46553.184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 8623
$ws.Range("I126").Value = 2019.8
$ws.Range("J126").Value = 12750
$ws.Range("K126").Value = 6059.4
$ws.Range("L126").Value = 38250
$ws.Range("M126").Value = -3589.4
$ws.Range("N126").Value = -43190

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2841855
$ws.Range("J107").Value = 5209449.5
$ws.Range("L107").Value = 15628348.5
$ws.Range("N107").Value = -15632188.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 4485.4
$ws.Range("J108").Value = 7100
$ws.Range("L108").Value = 21300
$ws.Range("N108").Value = -27060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 14566.5
$ws.Range("I110").Value = 5599.75
$ws.Range("J110").Value = 32500
$ws.Range("K110").Value = 16799.25
$ws.Range("L110").Value = 97500
$ws.Range("M110").Value = -12709.25
$ws.Range("N110").Value = -105680

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 2500
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 25001950
$ws.Range("J116").Value = 3250
$ws.Range("L116").Value = 9750
$ws.Range("N116").Value = -16634

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 4952.2
$ws.Range("I117").Value = 250
$ws.Range("J117").Value = 6127.75
$ws.Range("K117").Value = 750
$ws.Range("L117").Value = 18383.25
$ws.Range("M117").Value = 2692
$ws.Range("N117").Value = -25267.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 8024.5
$ws.Range("I118").Value = 699.3333
$ws.Range("K118").Value = 2097.9999
$ws.Range("M118").Value = -854.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1600.1852
$ws.Range("I132").Value = 1665.5714
$ws.Range("J132").Value = 1529.7693
$ws.Range("K132").Value = 14990.1426
$ws.Range("L132").Value = 13767.9237
$ws.Range("M132").Value = -12460.1426
$ws.Range("N132").Value = -18827.9237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4999.5
$ws.Range("I70").Value = 4999.5
$ws.Range("K70").Value = 4999.5
$ws.Range("M70").Value = -4729.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4999.5
$ws.Range("I73").Value = 4999.5
$ws.Range("K73").Value = 4999.5
$ws.Range("M73").Value = -4063.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 36748.5
$ws.Range("J92").Value = 36748.5
$ws.Range("L92").Value = 36748.5
$ws.Range("N92").Value = -40492.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1317.76
$ws.Range("I102").Value = 1272.8572
$ws.Range("K102").Value = 1272.8572
$ws.Range("M102").Value = 349.1428000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 993.5577
$ws.Range("I16").Value = 917.2
$ws.Range("K16").Value = 917.2
$ws.Range("M16").Value = -747.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 54
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 54
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 54
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -402

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2650.4167
$ws.Range("I2").Value = 2650.4167
$ws.Range("K2").Value = 2650.4167
$ws.Range("M2").Value = -2538.4167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1767.3334
$ws.Range("I81").Value = 1731.6
$ws.Range("J81").Value = 1812
$ws.Range("K81").Value = 3463.2
$ws.Range("L81").Value = 3624
$ws.Range("M81").Value = -2402.2
$ws.Range("N81").Value = -5746

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1767.3334
$ws.Range("I84").Value = 1731.6
$ws.Range("J84").Value = 1812
$ws.Range("K84").Value = 17316
$ws.Range("L84").Value = 18120
$ws.Range("M84").Value = -12012
$ws.Range("N84").Value = -28728

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 55559836
$ws.Range("I122").Value = 111113450
$ws.Range("K122").Value = 333340350
$ws.Range("M122").Value = -333337900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 59909
$ws.Range("J130").Value = 59909
$ws.Range("L130").Value = 59909
$ws.Range("N130").Value = -69949

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20195.834
$ws.Range("I136").Value = 764.9091
$ws.Range("K136").Value = 2294.7273
$ws.Range("M136").Value = 255.2727
